$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.912.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.14%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.573.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.87%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "549.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.83"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.94%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.00%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.05%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.51"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.50%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.029.36"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.48"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "61.819.62"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.08%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.28%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.580.47"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.72%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.58"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -3.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.54"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "337.42"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.01"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.491"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.55"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.29%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.04%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.24%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.10"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.12%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.24"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.87%  "

$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0832"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.11%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.34"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.83%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.78"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.57%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.04%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "19.12"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.40"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.14%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "326.63"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.47%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.02"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.899"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.92"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "37.48"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.61"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.605"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.95"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0544"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.50%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0964"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.25%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.45"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.68%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0238"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.86%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.059.34"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.41%  "
